$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) percentage columns
# Force text number-format on the touched cells so values such as
# '593.26' or '1.00' are stored verbatim as text, matching the source data
# (the sheet uses dotted/percent text, not real Excel numbers).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.697.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.612.10'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.26'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.94'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.548'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.607.71'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.126'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +9.15%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.48'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000185'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.077.44'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.539.96'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.609.07'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.16'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '363.56'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.66'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.95%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.01'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.35'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.74'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -8.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.742.04'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '575.97'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.93%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.93'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.86'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.67%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.71'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.32'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.32'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.84'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.55'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '155.51'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0288'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.98%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.82'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.50%  '
